$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: D3 loses its numeric value -> becomes a blank text cell ---
$ws.Range("D3").Value = "'"
$ws.Range("B2").Copy()
$ws.Range("D3").PasteSpecial(-4122)   # xlPasteFormats - restore D3's plain (unstyled) look

# --- Row 4: C4 corrected value 14990.33233432005 -> 0 ---
$ws.Range("C4").Value = 0

# --- Row 5: C5 corrected value 61189.18934749183 -> 3897.912874983024 ---
$ws.Range("C5").Value = 3897.912874983024

# --- Row 7: label "Other" -> "Biogas", value 129.0462505802004 -> 27.29972995275858 ---
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 27.29972995275858

# --- New row 8: "Other" row (split out of the old row 7 "Other" total) ---
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats - match the bold/bordered label style
$ws.Range("A8").Value = "Other"

$ws.Range("B8").Value = "'"
$ws.Range("C8").Value = "'"
$ws.Range("B2").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)   # xlPasteFormats - plain blank-cell look for B8/C8

$ws.Range("D8").Value = 109.4140670979836
